$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column E: Route 4 data.
# Header has no special style; data rows use a leading apostrophe so Excel
# applies the same quote-prefixed "text" style used by the sibling columns.
$ws.Range("E1").Value = "Route 4"
$ws.Range("E2").Value = "'BHM,EUS"
$ws.Range("E3").Value = "'EUS,BHM"
$ws.Range("E4").Value = "'0900,1300,1600"
$ws.Range("E5").Value = "'0700,0800,0900"
$ws.Range("E6").Value = "'0600,1400,1500"
$ws.Range("E7").Value = "'0900,1300,1600"
$ws.Range("E8").Value = "'0700,0800,0900"
$ws.Range("E9").Value = "'0600,1400,1500"

$ws.Range("E12").Select()
